$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.864.08"
$ws.Range("E2").Value = "  -1.44%  "
$ws.Range("D3").Value = "3.410.04"
$ws.Range("E3").Value = "  -1.14%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.29%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.411.57"
$ws.Range("E8").Value = "  -1.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.477"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.72%  "
$ws.Range("E10").Value = "  -1.56%  "
$ws.Range("E11").Value = "  +0.45%  "
$ws.Range("E12").Value = "  +0.86%  "
$ws.Range("D13").Value = "3.992.20"
$ws.Range("E13").Value = "  -1.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.29"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.86%  "
$ws.Range("E15").Value = "  +0.54%  "
$ws.Range("E16").Value = "  -1.77%  "
$ws.Range("D17").Value = "3.405.85"
$ws.Range("E17").Value = "  -1.35%  "
$ws.Range("D18").Value = "60.921.83"
$ws.Range("E18").Value = "  -1.41%  "
$ws.Range("E19").Value = "  -0.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "392.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.42%  "
$ws.Range("E23").Value = "  +0.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "72.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.72%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.997"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000122"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.45%  "
$ws.Range("D27").Value = "3.552.22"
$ws.Range("E27").Value = "  -1.05%  "
$ws.Range("E28").Value = "  +1.01%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.43"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.88%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.45"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.59%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.16"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.60%  "
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.89"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.66%  "
$ws.Range("E36").Value = "  -0.19%  "
$ws.Range("D37").Value = "3.436.64"
$ws.Range("E37").Value = "  -0.96%  "
$ws.Range("E38").Value = "  -2.38%  "
$ws.Range("E39").Value = "  -1.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "167.62"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0783"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.98"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.40%  "
$ws.Range("E43").Value = "  +0.32%  "
$ws.Range("E44").Value = "  +0.49%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.71"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.09%  "
$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "41.79"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.65%  "
$ws.Range("D48").Value = "2.592.09"
$ws.Range("E48").Value = "  -1.74%  "
$ws.Range("E49").Value = "  -4.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.96"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.93"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.96%  "
